# Daily update at 8 AM UTC
# Append the next day's row of data to the bottom of the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 40

$ws.Cells.Item($row, 1).Value = 45989
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 2).Value = 88
$ws.Cells.Item($row, 3).Value = 100
$ws.Cells.Item($row, 4).Value = 95
